$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to match the refreshed cryptos feed.
# NumberFormat="@" forces the assigned text to remain a text string
# (Excel would otherwise auto-coerce values like "1.002" or "2.880"
# into numbers). Resetting the Style back to "Normal" afterwards
# avoids leaving a residual text number-format on the cell.
$cellUpdates = @{
    "D2" = "25.704.59"
    "E2" = "  -5.56%  "
    "D3" = "1.812.43"
    "E3" = "  -4.85%  "
    "D4" = "1.002"
    "E4" = "  +0.11%  "
    "D5" = "276.76"
    "E5" = "  -9.61%  "
    "E6" = "  +0.10%  "
    "D7" = "0.5008"
    "E7" = "  -6.59%  "
    "D8" = "0.3494"
    "D9" = "44.15"
    "E9" = "  -3.67%  "
    "D10" = "0.06618"
    "E10" = "  -9.21%  "
    "D11" = "20.06"
    "E11" = "  -9.63%  "
    "D12" = "0.8417"
    "E12" = "  -7.05%  "
    "D13" = "0.07830"
    "E13" = "  -4.45%  "
    "D14" = "1.806.83"
    "E14" = "  +67.55%  "
    "E15" = "  -5.54%  "
    "D16" = "87.42"
    "E16" = "  -8.89%  "
    "D17" = "1.001"
    "E17" = "  +0.01%  "
    "E18" = "  -6.55%  "
    "E19" = "  +0.07%  "
    "D20" = "0.000007971"
    "E20" = "  -7.91%  "
    "D21" = "25.789.96"
    "E21" = "  -5.32%  "
    "D22" = "4.720"
    "E22" = "  -6.44%  "
    "D23" = "10.02"
    "E23" = "  -7.07%  "
    "D24" = "6.071"
    "E24" = "  -6.87%  "
    "D25" = "141.02"
    "E25" = "  -5.59%  "
    "D26" = "2.098"
    "E26" = "  -8.56%  "
    "D27" = "1.660"
    "E27" = "  -5.02%  "
    "D28" = "16.82"
    "E28" = "  -8.49%  "
    "D29" = "108.42"
    "E29" = "  -7.19%  "
    "D30" = "4.319"
    "E30" = "  -10.37%  "
    "D31" = "4.200"
    "E31" = "  -11.21%  "
    "D32" = "0.08767"
    "E32" = "  -4.92%  "
    "D33" = "0.04839"
    "E33" = "  -4.70%  "
    "D34" = "0.7386"
    "E34" = "  -11.03%  "
    "B35" = "ARBITRUM"
    "C35" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D35" = "1.131"
    "E35" = "  -7.16%  "
    "B36" = "HuobiToken"
    "C36" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D36" = "2.880"
    "E36" = "  -4.17%  "
    "E37" = "  +0.13%  "
    "D38" = "3.047"
    "E38" = "  -8.68%  "
    "D39" = "2.473"
    "E39" = "  -7.67%  "
    "D40" = "0.5324"
    "E40" = "  -8.97%  "
    "D41" = "0.01869"
    "E41" = "  -6.75%  "
    "D42" = "0.9696"
    "E42" = "  -9.93%  "
    "B43" = "FraxShare"
    "C43" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D43" = "6.230"
    "E43" = "  -5.91%  "
    "B44" = "Quant"
    "C44" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "D44" = "111.68"
    "E44" = "  -4.77%  "
    "D45" = "8.149"
    "E45" = "  -12.69%  "
    "D46" = "0.4685"
    "E46" = "  -7.58%  "
    "D47" = "1.000"
    "D48" = "0.1387"
    "E48" = "  -8.98%  "
    "D49" = "9.173"
    "E49" = "  -9.04%  "
    "E50" = "  -7.24%  "
    "B51" = "NEARProtocol"
    "C51" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D51" = "1.488"
    "E51" = "  -9.26%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
